$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- New data rows (39-43) to append to the TestData sheet ----
$rows = @(
    @{ num = 38; b = "validate specific user site data";
       f = "items[1].ID;items[1].VALUE;uris[1];etags[1]";
       g = '12090;Surgut;/mobile/custom/sdranalytics/users/1347/sites/12090;"1f-XDrwCbb23dQiQJkki2qr8wCx73Y"' },
    @{ num = 39; b = "validate reports menu cache";
       f = "items.REPORTID;items.REPORTHEADER_CUSTOMER;items.REPORTHEADER_RONUMBER;items.DEVICEDETAILS_TAG;items.CREATEDBY;uris;etags";
       g = '[dlcIaFJHr9oDyhlEnPxJ9qe5U2W1SOyEfqhT];[change Khushboo];[1234567];[abcdef];[1];[/mobile/custom/sdrmvp/reports/dlcIaFJHr9oDyhlEnPxJ9qe5U2W1SOyEfqhT/menus/dlcIaFJHr9oDyhlEnPxJ9qe5U2W1SOyEfqhT];["338-nehNJQHryyTFs/NJ5wlJ1X1GapU"]' },
    @{ num = 40; b = "validate report devices cache";
       f = "items.REPORTID;items.DEVICETYPE;items.PROCESSTYPE;items.CREATEDBY;uris;etags";
       g = '[dlcIaFJHr9oDyhlEnPxJ9qe5U2W1SOyEfqhT];[1];[340];[1];[/mobile/custom/sdrmvp/reports/dlcIaFJHr9oDyhlEnPxJ9qe5U2W1SOyEfqhT/devices/dlcIaFJHr9oDyhlEnPxJ9qe5U2W1SOyEfqhT];["80-rZanKi3FbEvjfwAHmlLn+xBgmxw"]' },
    @{ num = 41; b = "validate report headers cache";
       f = "items.ID;items.RO_NUMBER;items.FIELD_SERVICE_DIAGNOSTIC_ONLY;items.REPAIR_TYPE;items.REPORT_STATUS;items.EMPLOYEE_NAME;items.CUSTOMER_NAME;items.SERVICE_SITE;uris;etags";
       g = '[dlcIaFJHr9oDyhlEnPxJ9qe5U2W1SOyEfqhT];[1234567];[NO];[Field];[436];[Khushboo J];[change Khushboo];[12096];[/mobile/custom/sdrmvp/reports/dlcIaFJHr9oDyhlEnPxJ9qe5U2W1SOyEfqhT/headers/dlcIaFJHr9oDyhlEnPxJ9qe5U2W1SOyEfqhT];["226-6jVbyEDOOt7uhpitl2LT2S8ezuM"]' },
    @{ num = 42; b = "validate InprogressReports Cache";
       f = "items[1].ID;items[1].RO_NUMBER;items[1].FIELD_SERVICE_DIAGNOSTIC_ONLY;items[1].REPAIR_TYPE;items[1].REPORT_STATUS;items[1].EMPLOYEE_NAME;items[1].CUSTOMER_NAME;items[1].SITE_NAME;uris[1];etags[1]";
       g = 'PUWrgePHgrMT2zZxxL6vD5h3mZfBBJfhJ2uu;R3256;NO;Depot;436;Parul Gupta;Khushboo;ABC1;/mobile/custom/sdranalytics/users/1/inprogressreports/PUWrgePHgrMT2zZxxL6vD5h3mZfBBJfhJ2uu;"257-RhuXjhxj6FQctD5qBYVQ0KQ6IEI"' }
)

# Template cells (already-styled, bordered rows) used to replicate formatting
$srcAE = $ws.Range("A37:E37")
$srcFG = $ws.Range("F37:G37")

$startRow = 39
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    # Copy formatting from an existing formatted row onto the new row
    $srcAE.Copy()
    $ws.Range("A$r`:E$r").PasteSpecial(-4122)
    $srcFG.Copy()
    $ws.Range("F$r`:G$r").PasteSpecial(-4122)

    $ws.Cells.Item($r, 1).Value = $data.num
    $ws.Cells.Item($r, 2).Value = $data.b
    $ws.Cells.Item($r, 3).Value = "Emerson"
    $ws.Cells.Item($r, 4).Value = "GET"
    $ws.Cells.Item($r, 5).Value = "{}"
    $ws.Cells.Item($r, 6).Value = $data.f
    $ws.Cells.Item($r, 7).Value = $data.g
}

$excel.CutCopyMode = 0

# ---- Update sheet view state to reflect the new active selection ----
$ws.Application.ActiveWindow.ScrollRow = 40
$ws.Range("G43").Select()
